$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 129
$ws.Range("H129").Value = 906.4706
$ws.Range("J129").Value = 1003.6842
$ws.Range("L129").Value = 3011.0526
$ws.Range("N129").Value = -13011.0526

# Row 135
$ws.Range("H135").Value = 22088.596
$ws.Range("I135").Value = 22502.826
$ws.Range("K135").Value = 202525.434
$ws.Range("M135").Value = -199990.434

# Row 137
$ws.Range("H137").Value = 1925007.8
$ws.Range("I137").Value = 3126246.8
$ws.Range("J137").Value = 3025.6
$ws.Range("K137").Value = 9378740.399999999
$ws.Range("L137").Value = 9076.799999999999
$ws.Range("M137").Value = -9376190.399999999
$ws.Range("N137").Value = -14176.8

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1297.6923
$ws.Range("I2").Value = 1224.5454
$ws.Range("J2").Value = 1700
$ws.Range("K2").Value = 1224.5454
$ws.Range("L2").Value = 1700
$ws.Range("M2").Value = -1111.5454
$ws.Range("N2").Value = -1926

# Row 11
$ws.Range("H11").Value = 7000
$ws.Range("I11").Value = 10000
$ws.Range("J11").Value = 4000
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 4000
$ws.Range("M11").Value = -9856
$ws.Range("N11").Value = -4288

# Row 61
$ws.Range("H61").Value = 22267798
$ws.Range("I61").Value = 27806534
$ws.Range("J61").Value = 112853.336
$ws.Range("K61").Value = 27806534
$ws.Range("L61").Value = 112853.336
$ws.Range("M61").Value = -27806322
$ws.Range("N61").Value = -113277.336

# Row 74
$ws.Range("H74").Value = 4847990
$ws.Range("I74").Value = 6784895
$ws.Range("K74").Value = 6784895
$ws.Range("M74").Value = -6784021

# Row 77
$ws.Range("H77").Value = 4847990
$ws.Range("I77").Value = 6784895
$ws.Range("K77").Value = 33924475
$ws.Range("M77").Value = -33920107

# Row 116
$ws.Range("H116").Value = 1297.6923
$ws.Range("I116").Value = 1224.5454
$ws.Range("J116").Value = 1700
$ws.Range("K116").Value = 1224.5454
$ws.Range("L116").Value = 1700
$ws.Range("M116").Value = 1069.4546
$ws.Range("N116").Value = -6288

# Row 132
$ws.Range("H132").Value = 48599.37
$ws.Range("I132").Value = 38637.668
$ws.Range("J132").Value = 65409.75
$ws.Range("K132").Value = 115913.004
$ws.Range("L132").Value = 196229.25
$ws.Range("M132").Value = -113383.004
$ws.Range("N132").Value = -201289.25

# Row 135
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

# Row 136
$ws.Range("H136").Value = 22267798
$ws.Range("I136").Value = 27806534
$ws.Range("J136").Value = 112853.336
$ws.Range("K136").Value = 83419602
$ws.Range("L136").Value = 338560.008
$ws.Range("M136").Value = -83417052
$ws.Range("N136").Value = -343660.008

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1297.6923
$ws.Range("I3").Value = 1224.5454
$ws.Range("J3").Value = 1700
$ws.Range("K3").Value = 1224.5454
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = -1110.5454
$ws.Range("N3").Value = -1928

# Row 134
$ws.Range("H134").Value = 2263.9707
$ws.Range("I134").Value = 1594.8518
$ws.Range("K134").Value = 4784.555399999999
$ws.Range("M134").Value = -2249.555399999999

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 58824090
$ws.Range("I22").Value = 111111430
$ws.Range("K22").Value = 111111430
$ws.Range("M22").Value = -111111080

# Row 31
$ws.Range("H31").Value = 1647.1025
$ws.Range("I31").Value = 1175.9032
$ws.Range("J31").Value = 3473
$ws.Range("K31").Value = 1175.9032
$ws.Range("L31").Value = 3473
$ws.Range("M31").Value = -880.9032
$ws.Range("N31").Value = -4063

# Row 34
$ws.Range("H34").Value = 1647.1025
$ws.Range("I34").Value = 1175.9032
$ws.Range("J34").Value = 3473
$ws.Range("K34").Value = 1175.9032
$ws.Range("L34").Value = 3473
$ws.Range("M34").Value = -973.9032
$ws.Range("N34").Value = -3877

# Row 58
$ws.Range("H58").Value = 17545326
$ws.Range("I58").Value = 25642052
$ws.Range("J58").Value = 2418.2222
$ws.Range("K58").Value = 25642052
$ws.Range("L58").Value = 2418.2222
$ws.Range("M58").Value = -25641849
$ws.Range("N58").Value = -2824.2222

# Row 99
$ws.Range("H99").Value = 2692.1304
$ws.Range("I99").Value = 2437.1667
$ws.Range("J99").Value = 3610
$ws.Range("K99").Value = 2437.1667
$ws.Range("L99").Value = 3610
$ws.Range("M99").Value = -939.1667000000002
$ws.Range("N99").Value = -6606

# Row 126
$ws.Range("H126").Value = 2692.1304
$ws.Range("I126").Value = 2437.1667
$ws.Range("J126").Value = 3610
$ws.Range("K126").Value = 7311.500100000001
$ws.Range("L126").Value = 10830
$ws.Range("M126").Value = -4841.500100000001
$ws.Range("N126").Value = -15770

# Row 132
$ws.Range("H132").Value = 31041.771
$ws.Range("I132").Value = 2157.739
$ws.Range("J132").Value = 86402.836
$ws.Range("K132").Value = 6473.217000000001
$ws.Range("L132").Value = 259208.508
$ws.Range("M132").Value = -3943.217000000001
$ws.Range("N132").Value = -264268.508

# Row 134
$ws.Range("H134").Value = 22773.793
$ws.Range("I134").Value = 1496.9773
$ws.Range("J134").Value = 126793.78
$ws.Range("K134").Value = 4490.9319
$ws.Range("L134").Value = 380381.34
$ws.Range("M134").Value = -1955.9319
$ws.Range("N134").Value = -385451.34

# Row 136
$ws.Range("H136").Value = 17545326
$ws.Range("I136").Value = 25642052
$ws.Range("J136").Value = 2418.2222
$ws.Range("K136").Value = 76926156
$ws.Range("L136").Value = 7254.6666
$ws.Range("M136").Value = -76923606
$ws.Range("N136").Value = -12354.6666

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 19392.326
$ws.Range("I131").Value = 878
$ws.Range("J131").Value = 21496.227
$ws.Range("K131").Value = 2634
$ws.Range("L131").Value = 64488.681
$ws.Range("M131").Value = 2406
$ws.Range("N131").Value = -74568.681

# Row 136
$ws.Range("H136").Value = 1669.6
$ws.Range("I136").Value = 783.3333
$ws.Range("J136").Value = 2999
$ws.Range("K136").Value = 2349.9999
$ws.Range("L136").Value = 8997
$ws.Range("M136").Value = 2750.0001
$ws.Range("N136").Value = -19197

# Row 137
$ws.Range("H137").Value = 23670.555
$ws.Range("I137").Value = 982.3333
$ws.Range("J137").Value = 41821.133
$ws.Range("K137").Value = 2946.9999
$ws.Range("L137").Value = 125463.399
$ws.Range("M137").Value = 2153.0001
$ws.Range("N137").Value = -135663.399

# Row 140
$ws.Range("H140").Value = 2189.8901
$ws.Range("I140").Value = 2027.2727
$ws.Range("K140").Value = 6081.8181
$ws.Range("M140").Value = -901.8181000000004

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 10000
$ws.Range("J5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("N5").Value = -10224

# Row 132
$ws.Range("H132").Value = 70721.34
$ws.Range("I132").Value = 63538.625
$ws.Range("J132").Value = 79561.62
$ws.Range("K132").Value = 190615.875
$ws.Range("L132").Value = 238684.86
$ws.Range("M132").Value = -188085.875
$ws.Range("N132").Value = -243744.86

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 894.65
$ws.Range("I22").Value = 487.625
$ws.Range("J22").Value = 1166
$ws.Range("K22").Value = 487.625
$ws.Range("L22").Value = 1166
$ws.Range("M22").Value = -192.625
$ws.Range("N22").Value = -1756

# Row 27
$ws.Range("H27").Value = 894.65
$ws.Range("I27").Value = 487.625
$ws.Range("J27").Value = 1166
$ws.Range("K27").Value = 487.625
$ws.Range("L27").Value = 1166
$ws.Range("M27").Value = -380.625
$ws.Range("N27").Value = -1380

# Row 122
$ws.Range("H122").Value = 2613.12
$ws.Range("I122").Value = 2574
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 7722
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -5272
$ws.Range("N122").Value = -13600

# Row 132
$ws.Range("H132").Value = 18544.88
$ws.Range("I132").Value = 1320.7609
$ws.Range("J132").Value = 79491.766
$ws.Range("K132").Value = 3962.2827
$ws.Range("L132").Value = 238475.298
$ws.Range("M132").Value = -1432.2827
$ws.Range("N132").Value = -243535.298

# Row 136
$ws.Range("H136").Value = 75321.516
$ws.Range("I136").Value = 56610.5
$ws.Range("J136").Value = 112743.555
$ws.Range("K136").Value = 169831.5
$ws.Range("L136").Value = 338230.665
$ws.Range("M136").Value = -167281.5
$ws.Range("N136").Value = -343330.665

$ws = $wb.Worksheets.Item("WVR")
# Row 68
$ws.Range("H68").Value = 44092
$ws.Range("J68").Value = 44092
$ws.Range("L68").Value = 44092
$ws.Range("N68").Value = -45714

# Row 71
$ws.Range("H71").Value = 44092
$ws.Range("J71").Value = 44092
$ws.Range("L71").Value = 132276
$ws.Range("N71").Value = -140388

# Row 132
$ws.Range("H132").Value = 45490.156
$ws.Range("I132").Value = 29228.514
$ws.Range("J132").Value = 102405.9
$ws.Range("K132").Value = 87685.542
$ws.Range("L132").Value = 307217.7
$ws.Range("M132").Value = -85155.542
$ws.Range("N132").Value = -312277.7

# Row 136
$ws.Range("H136").Value = 37427.656
$ws.Range("I136").Value = 22206.766
$ws.Range("J136").Value = 126850.375
$ws.Range("K136").Value = 66620.298
$ws.Range("L136").Value = 380551.125
$ws.Range("M136").Value = -64070.298
$ws.Range("N136").Value = -385651.125

